$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1941.2142
$ws.Range("J19").Value = 1542.8889
$ws.Range("L19").Value = 1542.8889
$ws.Range("N19").Value = -1892.8889

# Row 103
$ws.Range("H103").Value = 30050002
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 107
$ws.Range("H107").Value = 6166.8
$ws.Range("I107").Value = 7583.5
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 7583.5
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = -5663.5
$ws.Range("N107").Value = -4340

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 343.33334
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 186.66667
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 186.66667
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = -418.66667

# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()

# Row 88
$ws.Range("H88").Value = 2425.7917
$ws.Range("I88").Value = 2232
$ws.Range("J88").Value = 2813.375
$ws.Range("K88").Value = 2232
$ws.Range("L88").Value = 2813.375
$ws.Range("M88").Value = -1826
$ws.Range("N88").Value = -3625.375

# Row 91
$ws.Range("H91").Value = 2425.7917
$ws.Range("I91").Value = 2232
$ws.Range("J91").Value = 2813.375
$ws.Range("K91").Value = 2232
$ws.Range("L91").Value = 2813.375
$ws.Range("M91").Value = -828
$ws.Range("N91").Value = -5621.375

# Row 132
$ws.Range("H132").Value = 6946397.5
$ws.Range("I132").Value = 9616804
$ws.Range("J132").Value = 3341.2
$ws.Range("K132").Value = 28850412
$ws.Range("L132").Value = 10023.6
$ws.Range("M132").Value = -28847882
$ws.Range("N132").Value = -15083.6

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()

# Row 86
$ws.Range("H86").Value = 20835644
$ws.Range("I86").Value = 2430.8333
$ws.Range("J86").Value = 41668856
$ws.Range("K86").Value = 2430.8333
$ws.Range("L86").Value = 41668856
$ws.Range("M86").Value = -1307.8333
$ws.Range("N86").Value = -41671102

# Row 89
$ws.Range("H89").Value = 20835644
$ws.Range("I89").Value = 2430.8333
$ws.Range("J89").Value = 41668856
$ws.Range("K89").Value = 12154.1665
$ws.Range("L89").Value = 208344280
$ws.Range("M89").Value = -6538.166499999999
$ws.Range("N89").Value = -208355512

# Row 107
$ws.Range("H107").Value = 1116.7142
$ws.Range("I107").Value = 921.2727
$ws.Range("J107").Value = 1833.3334
$ws.Range("K107").Value = 921.2727
$ws.Range("L107").Value = 1833.3334
$ws.Range("M107").Value = 998.7273
$ws.Range("N107").Value = -5673.3334

# Row 134
$ws.Range("H134").Value = 3581.3958
$ws.Range("I134").Value = 2741.4167
$ws.Range("J134").Value = 6101.3335
$ws.Range("K134").Value = 8224.250100000001
$ws.Range("L134").Value = 18304.0005
$ws.Range("M134").Value = -5689.250100000001
$ws.Range("N134").Value = -23374.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6176057
$ws.Range("I31").Value = 3338.28
$ws.Range("J31").Value = 83335040
$ws.Range("K31").Value = 3338.28
$ws.Range("L31").Value = 83335040
$ws.Range("M31").Value = -3043.28
$ws.Range("N31").Value = -83335630

# Row 34
$ws.Range("H34").Value = 6176057
$ws.Range("I34").Value = 3338.28
$ws.Range("J34").Value = 83335040
$ws.Range("K34").Value = 3338.28
$ws.Range("L34").Value = 83335040
$ws.Range("M34").Value = -3136.28
$ws.Range("N34").Value = -83335444

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("L92").ClearContents()

# Row 107
$ws.Range("H107").Value = 914.9474
$ws.Range("I107").Value = 604.8
$ws.Range("J107").Value = 2078
$ws.Range("K107").Value = 604.8
$ws.Range("L107").Value = 2078
$ws.Range("M107").Value = 1315.2
$ws.Range("N107").Value = -5918

# Row 132
$ws.Range("H132").Value = 27780320
$ws.Range("I132").Value = 41668120
$ws.Range("J132").Value = 4718.6665
$ws.Range("K132").Value = 125004360
$ws.Range("L132").Value = 14155.9995
$ws.Range("M132").Value = -125001830
$ws.Range("N132").Value = -19215.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 2380.6365
$ws.Range("I80").Value = 2328
$ws.Range("K80").Value = 6984
$ws.Range("M80").Value = -6048

# Row 83
$ws.Range("H83").Value = 2380.6365
$ws.Range("I83").Value = 2328
$ws.Range("K83").Value = 20952
$ws.Range("M83").Value = -16272

# Row 93
$ws.Range("H93").Value = 4900
$ws.Range("J93").Value = 4900
$ws.Range("L93").Value = 14700
$ws.Range("N93").Value = -18444

# Row 94
$ws.Range("H94").Value = 4380
$ws.Range("I94").Value = 5000
$ws.Range("J94").Value = 4311.1113
$ws.Range("K94").Value = 15000
$ws.Range("L94").Value = 12933.3339
$ws.Range("M94").Value = -14324
$ws.Range("N94").Value = -14285.3339

# Row 95
$ws.Range("H95").Value = 12660
$ws.Range("I95").Value = 10000
$ws.Range("J95").Value = 13990
$ws.Range("K95").Value = 30000
$ws.Range("L95").Value = 41970
$ws.Range("N95").Value = -46088
$ws.Range("M95").Value = -27941

# Row 97
$ws.Range("H97").Value = 9406.929
$ws.Range("I97").Value = 34401
$ws.Range("J97").Value = 2590.3635
$ws.Range("K97").Value = 103203
$ws.Range("L97").Value = 7771.0905
$ws.Range("M97").Value = -102707
$ws.Range("N97").Value = -8763.0905

# Row 99
$ws.Range("H99").Value = 2186.182
$ws.Range("I99").Value = 774
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2322
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -76
$ws.Range("N99").Value = -11992

# Row 100
$ws.Range("H100").Value = 4605.3335
$ws.Range("J100").Value = 4605.3335
$ws.Range("L100").Value = 13816.0005
$ws.Range("N100").Value = -15438.0005

# Row 101
$ws.Range("H101").Value = 10575
$ws.Range("J101").Value = 10575
$ws.Range("L101").Value = 31725
$ws.Range("N101").Value = -36593

# Row 102
$ws.Range("H102").Value = 13555.556
$ws.Range("J102").Value = 13555.556
$ws.Range("L102").Value = 40666.66800000001
$ws.Range("N102").Value = -45534.66800000001

# Row 104
$ws.Range("H104").Value = 3749.375
$ws.Range("J104").Value = 1499.3334
$ws.Range("L104").Value = 4498.0002
$ws.Range("N104").Value = -9740.0002

# Row 106
$ws.Range("H106").Value = 10360
$ws.Range("J106").Value = 10360
$ws.Range("L106").Value = 31080
$ws.Range("N106").Value = -32972

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("N117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("M117").ClearContents()

# Row 118
$ws.Range("H118").Value = 1517.2
$ws.Range("I118").Value = 290
$ws.Range("J118").Value = 1653.5555
$ws.Range("K118").Value = 870
$ws.Range("L118").Value = 4960.666499999999
$ws.Range("M118").Value = 373
$ws.Range("N118").Value = -7446.666499999999

# Row 121
$ws.Range("H121").Value = 931.2353000000001
$ws.Range("I121").Value = 277.14285
$ws.Range("J121").Value = 1389.1
$ws.Range("K121").Value = 831.4285500000001
$ws.Range("L121").Value = 4167.299999999999
$ws.Range("M121").Value = 478.5714499999999
$ws.Range("N121").Value = -6787.299999999999

# Row 129
$ws.Range("H129").Value = 2253.5833
$ws.Range("I129").Value = 1180
$ws.Range("J129").Value = 2897.7334
$ws.Range("K129").Value = 3540
$ws.Range("L129").Value = 8693.200199999999
$ws.Range("M129").Value = 1460
$ws.Range("N129").Value = -18693.2002

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 10104359
$ws.Range("I80").Value = 14495310
$ws.Range("J80").Value = 5171.2
$ws.Range("K80").Value = 14495310
$ws.Range("L80").Value = 5171.2
$ws.Range("M80").Value = -14494312
$ws.Range("N80").Value = -7167.2

# Row 83
$ws.Range("H83").Value = 10104359
$ws.Range("I83").Value = 14495310
$ws.Range("J83").Value = 5171.2
$ws.Range("K83").Value = 72476550
$ws.Range("L83").Value = 25856
$ws.Range("M83").Value = -72471558
$ws.Range("N83").Value = -35840

# Row 122
$ws.Range("H122").Value = 3176521.5
$ws.Range("I122").Value = 4167966.5
$ws.Range("J122").Value = 3897.4
$ws.Range("K122").Value = 12503899.5
$ws.Range("L122").Value = 11692.2
$ws.Range("M122").Value = -12501449.5
$ws.Range("N122").Value = -16592.2

# Row 132
$ws.Range("H132").Value = 3577.182
$ws.Range("I132").Value = 2287.3809
$ws.Range("J132").Value = 5834.3335
$ws.Range("K132").Value = 6862.1427
$ws.Range("L132").Value = 17503.0005
$ws.Range("M132").Value = -4332.1427
$ws.Range("N132").Value = -22563.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1254.4286
$ws.Range("I16").Value = 1327.8462
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 1327.8462
$ws.Range("L16").Value = 300
$ws.Range("M16").Value = -1157.8462
$ws.Range("N16").Value = -640

# Row 81
$ws.Range("H81").Value = 35027.4
$ws.Range("J81").Value = 35027.4
$ws.Range("L81").Value = 35027.4
$ws.Range("N81").Value = -37023.4

# Row 84
$ws.Range("H84").Value = 35027.4
$ws.Range("J84").Value = 35027.4
$ws.Range("L84").Value = 105082.2
$ws.Range("N84").Value = -115066.2

# Row 122
$ws.Range("H122").Value = 6450.316
$ws.Range("I122").Value = 7335.6
$ws.Range("J122").Value = 5466.6665
$ws.Range("K122").Value = 22006.8
$ws.Range("L122").Value = 16399.9995
$ws.Range("M122").Value = -19556.8
$ws.Range("N122").Value = -21299.9995

# Row 136
$ws.Range("H136").Value = 9808221
$ws.Range("I136").Value = 12822377
$ws.Range("J136").Value = 12217.083
$ws.Range("K136").Value = 38467131
$ws.Range("L136").Value = 36651.249
$ws.Range("M136").Value = -38464581
$ws.Range("N136").Value = -41751.249
